$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation was recorded for 2023-06-16 (volume 80) ahead of the
# existing 2023-04-26 row. Insert a fresh row at 69 (pushing the old
# rows 69-85 down to 70-86) and seed it with a copy of what is now row 70
# (the former row 69), then overwrite the two cells that actually differ
# for the new record (Fecha / Volumen).
$ws.Rows("69:69").Insert()

$ws.Range("A70:T70").Copy() | Out-Null
$ws.Range("A69:T69").PasteSpecial() | Out-Null

$ws.Range("D69").Value = 45093
$ws.Range("M69").Value = 80
